$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = 26
$row2 = 27

# Columns whose values actually differ between row 26 and row 27 (per the
# source diff). Columns not listed here are left completely untouched so
# that empty/blank cells and text-like values (dates, times) are not
# re-typed by Excel's automatic value parsing.
$numericCols = @("A", "B", "E", "Q", "R", "S")
$textCols = @("D", "F", "G", "H", "I", "P", "Z", "AB", "AW", "AX")

foreach ($col in $numericCols) {
    $addr1 = $col + $row1
    $addr2 = $col + $row2
    $v1 = $ws.Range($addr1).Value()
    $v2 = $ws.Range($addr2).Value()
    $ws.Range($addr1).Value = $v2
    $ws.Range($addr2).Value = $v1
}

foreach ($col in $textCols) {
    $addr1 = $col + $row1
    $addr2 = $col + $row2
    $v1 = $ws.Range($addr1).Value()
    $v2 = $ws.Range($addr2).Value()

    # Force text format so Excel doesn't reinterpret numeric/date/time-like
    # strings (e.g. "10", "00:00", "19:51") as a different cell type.
    $ws.Range($addr1).NumberFormat = "@"
    $ws.Range($addr2).NumberFormat = "@"

    $ws.Range($addr1).Value = [string]$v2
    $ws.Range($addr2).Value = [string]$v1
}
